$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the existing "Rules" block (rows 33-41) down by two rows so that
#    new rows can be inserted above it (new row 33 block + blank separator).
# ---------------------------------------------------------------------------
$ws.Rows("33:34").Insert()

# ---------------------------------------------------------------------------
# 2. New row 33 : SARIF1011 / ReferToFinalSchema block
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "SARIF1011"
$ws.Range("B33").Value = "ReferToFinalSchema"
$ws.Range("C33").Value = "error"
$ws.Range("D33").Value = "Schema"
$ws.Range("E33").Value = "SchemaPropertyMustReferToFinalSchema"
$ws.Range("F16").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value = "DONE"

# F32 : leftover formatting-only cell (bold white font, no fill), no value.
$ws.Range("F16").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$ws.Range("F32").Font.ThemeColor = 14
$ws.Range("F32").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 3. Append the two extra rows for the existing "ProduceEnrichedSarif" block
#    (now at rows 41-45 after the shift above).
# ---------------------------------------------------------------------------
$ws.Range("D41").Value = "Content"

$ws.Range("E44").Value = "ProvideHelpUris"
$ws.Range("F5").Copy()
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("F44").Value = "TODO"

$ws.Range("E45").Value = "EmbedFileContent"
$ws.Range("F5").Copy()
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("F45").Value = "TODO"

# ---------------------------------------------------------------------------
# 4. New row 47-48 : SARIF2004 / ReduceFileSize block
# ---------------------------------------------------------------------------
$ws.Range("A47").Value = "SARIF2004"
$ws.Range("B47").Value = "ReduceFileSize"
$ws.Range("C47").Value = "warning"
$ws.Range("D47").Value = "Tool"
$ws.Range("E47").Value = "EliminateLocationOnlyArtifacts"
$ws.Range("F5").Copy()
$ws.Range("F47").PasteSpecial(-4122)
$ws.Range("F47").Value = "TODO"

$ws.Range("E48").Value = "DoNotIncludeExtraIndexedObjectProperties"
$ws.Range("F5").Copy()
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("F48").Value = "TODO"

# ---------------------------------------------------------------------------
# 5. New row 50-52 : SARIF2005 / ProvideHelpfulToolInformation block
# ---------------------------------------------------------------------------
$ws.Range("A50").Value = "SARIF2005"
$ws.Range("B50").Value = "ProvideHelpfulToolInformation"
$ws.Range("C50").Value = "warning"
$ws.Range("D50").Value = "Size"
$ws.Range("E50").Value = "ProvideConciseToolName"
$ws.Range("F5").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F50").Value = "TODO"

$ws.Range("E51").Value = "ProvideToolVersion"
$ws.Range("F5").Copy()
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F51").Value = "TODO"

$ws.Range("E52").Value = "UseNumericToolVersions"
$ws.Range("F5").Copy()
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("F52").Value = "TODO"

# ---------------------------------------------------------------------------
# 6. New row 54 : SARI2006 / ProvideUsableUris block
# ---------------------------------------------------------------------------
$ws.Range("A54").Value = "SARI2006"
$ws.Range("B54").Value = "ProvideUsableUris"
$ws.Range("C54").Value = "warning"
$ws.Range("D54").Value = "Uris"
$ws.Range("E54").Value = "UrisShouldBeReachable"
$ws.Range("F5").Copy()
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F54").Value = "TODO"

# ---------------------------------------------------------------------------
# 7. New row 56 : SARIF2007 / ExpressPathsRelativeToReplRoot block
# ---------------------------------------------------------------------------
$ws.Range("A56").Value = "SARIF2007"
$ws.Range("B56").Value = "ExpressPathsRelativeToReplRoot"
$ws.Range("C56").Value = "warning"
$ws.Range("E56").Value = "(name?)"
$ws.Range("F5").Copy()
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("F56").Value = "TODO"

# ---------------------------------------------------------------------------
# 8. New row 58 : SARIF2008 / ProvideSchema block
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "SARIF2008"
$ws.Range("B58").Value = "ProvideSchema"
$ws.Range("C58").Value = "warning"
$ws.Range("D58").Value = "Schema"
$ws.Range("E58").Value = "SchemaPropertyShouldBePresent"
$ws.Range("F16").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F58").Value = "DONE"

# ---------------------------------------------------------------------------
# 9. Column A width shrinks now that the long rule names moved to column B.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 11

# ---------------------------------------------------------------------------
# 10. Selection / view cosmetics.
# ---------------------------------------------------------------------------
$ws.Range("H2").Select()
